$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(4, 1).Value = "ryyy"
$ws.Cells.Item(4, 2).Value = "cameron"
$ws.Cells.Item(4, 3).Value = "testing x"
$ws.Cells.Item(4, 4).Value = "2025-09-27 00:41:48"
